$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 9165
$ws.Cells.Item(28, 9).Value = 3772.647
$ws.Cells.Item(28, 11).Value = 3772.647
$ws.Cells.Item(28, 13).Value = -3287.647
$ws.Cells.Item(32, 8).Value = 873
$ws.Cells.Item(32, 10).Value = 873
$ws.Cells.Item(32, 12).Value = 873
$ws.Cells.Item(32, 14).Value = -1525
$ws.Cells.Item(38, 8).Value = 815.7
$ws.Cells.Item(38, 9).Value = 154.33333
$ws.Cells.Item(38, 10).Value = 2799.8
$ws.Cells.Item(38, 11).Value = 462.99999
$ws.Cells.Item(38, 12).Value = 8399.400000000001
$ws.Cells.Item(38, 13).Value = -90.99998999999997
$ws.Cells.Item(38, 14).Value = -9143.400000000001
$ws.Cells.Item(51, 8).Value = 3833.3333
$ws.Cells.Item(51, 10).Value = 3833.3333
$ws.Cells.Item(51, 12).Value = 3833.3333
$ws.Cells.Item(51, 14).Value = -4801.3333
$ws.Cells.Item(58, 8).Value = 1379.4615
$ws.Cells.Item(58, 9).Value = 61.857143
$ws.Cells.Item(58, 10).Value = 2916.6667
$ws.Cells.Item(58, 11).Value = 185.571429
$ws.Cells.Item(58, 12).Value = 8750.000100000001
$ws.Cells.Item(58, 13).Value = -35.57142899999999
$ws.Cells.Item(58, 14).Value = -9050.000100000001
$ws.Cells.Item(69, 8).Value = 6980.909
$ws.Cells.Item(69, 9).Value = 6800
$ws.Cells.Item(69, 11).Value = 20400
$ws.Cells.Item(69, 13).Value = -19526
$ws.Cells.Item(72, 8).Value = 6980.909
$ws.Cells.Item(72, 9).Value = 6800
$ws.Cells.Item(72, 11).Value = 61200
$ws.Cells.Item(72, 13).Value = -56832
$ws.Cells.Item(76, 8).Value = 4749.25
$ws.Cells.Item(76, 10).Value = 4000
$ws.Cells.Item(76, 12).Value = 4000
$ws.Cells.Item(76, 14).Value = -4630
$ws.Cells.Item(79, 8).Value = 4749.25
$ws.Cells.Item(79, 10).Value = 4000
$ws.Cells.Item(79, 12).Value = 4000
$ws.Cells.Item(79, 14).Value = -6184
$ws.Cells.Item(121, 8).Value = 730.5714
$ws.Cells.Item(121, 10).Value = 730.5714
$ws.Cells.Item(121, 12).Value = 2191.7142
$ws.Cells.Item(121, 14).Value = -5685.7142
$ws.Cells.Item(129, 8).Value = 2476
$ws.Cells.Item(129, 10).Value = 3000
$ws.Cells.Item(129, 12).Value = 9000
$ws.Cells.Item(129, 14).Value = -19000
$ws.Cells.Item(131, 8).Value = 1394
$ws.Cells.Item(131, 9).Value = 1041.5
$ws.Cells.Item(131, 11).Value = 3124.5
$ws.Cells.Item(131, 13).Value = 1915.5
$ws.Cells.Item(138, 8).Value = 3273.2778
$ws.Cells.Item(138, 10).Value = 3527.9333
$ws.Cells.Item(138, 12).Value = 10583.7999
$ws.Cells.Item(138, 14).Value = -20863.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2955.9443
$ws.Cells.Item(2, 9).Value = 2530.4614
$ws.Cells.Item(2, 10).Value = 4062.2
$ws.Cells.Item(2, 11).Value = 2530.4614
$ws.Cells.Item(2, 12).Value = 4062.2
$ws.Cells.Item(2, 13).Value = -2417.4614
$ws.Cells.Item(2, 14).Value = -4288.2
$ws.Cells.Item(32, 8).Value = 13241.625
$ws.Cells.Item(32, 9).Value = 8528.191999999999
$ws.Cells.Item(32, 11).Value = 8528.191999999999
$ws.Cells.Item(32, 13).Value = -8241.191999999999
$ws.Cells.Item(45, 8).Value = 4783.1665
$ws.Cells.Item(45, 9).Value = 2500
$ws.Cells.Item(45, 11).Value = 2500
$ws.Cells.Item(45, 13).Value = -2123
$ws.Cells.Item(63, 8).Value = 2928.4285
$ws.Cells.Item(66, 8).Value = 2928.4285
$ws.Cells.Item(116, 8).Value = 2955.9443
$ws.Cells.Item(116, 9).Value = 2530.4614
$ws.Cells.Item(116, 10).Value = 4062.2
$ws.Cells.Item(116, 11).Value = 2530.4614
$ws.Cells.Item(116, 12).Value = 4062.2
$ws.Cells.Item(116, 13).Value = -236.4614000000001
$ws.Cells.Item(116, 14).Value = -8650.200000000001
$ws.Cells.Item(122, 8).Value = 1726.6666
$ws.Cells.Item(122, 9).Value = 1742.5
$ws.Cells.Item(122, 10).Value = 1600
$ws.Cells.Item(122, 11).Value = 5227.5
$ws.Cells.Item(122, 12).Value = 4800
$ws.Cells.Item(122, 13).Value = -2777.5
$ws.Cells.Item(122, 14).Value = -9700
$ws.Cells.Item(132, 8).Value = 1606.5
$ws.Cells.Item(132, 10).Value = 2599.8
$ws.Cells.Item(132, 12).Value = 7799.400000000001
$ws.Cells.Item(132, 14).Value = -12859.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2955.9443
$ws.Cells.Item(3, 9).Value = 2530.4614
$ws.Cells.Item(3, 10).Value = 4062.2
$ws.Cells.Item(3, 11).Value = 2530.4614
$ws.Cells.Item(3, 12).Value = 4062.2
$ws.Cells.Item(3, 13).Value = -2416.4614
$ws.Cells.Item(3, 14).Value = -4290.2
$ws.Cells.Item(86, 8).Value = 8999.799999999999
$ws.Cells.Item(86, 9).Value = 9000
$ws.Cells.Item(86, 11).Value = 9000
$ws.Cells.Item(86, 13).Value = -7877
$ws.Cells.Item(89, 8).Value = 8999.799999999999
$ws.Cells.Item(89, 9).Value = 9000
$ws.Cells.Item(89, 11).Value = 45000
$ws.Cells.Item(89, 13).Value = -39384
$ws.Cells.Item(94, 8).Value = 1542
$ws.Cells.Item(94, 9).Value = 1542
$ws.Cells.Item(94, 11).Value = 1542
$ws.Cells.Item(94, 13).Value = -1091
$ws.Cells.Item(107, 8).Value = 2893.6052
$ws.Cells.Item(107, 9).Value = 1365.9231
$ws.Cells.Item(107, 10).Value = 6203.5835
$ws.Cells.Item(107, 11).Value = 1365.9231
$ws.Cells.Item(107, 12).Value = 6203.5835
$ws.Cells.Item(107, 13).Value = 554.0769
$ws.Cells.Item(107, 14).Value = -10043.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 3336133
$ws.Cells.Item(3, 9).Value = 5001750
$ws.Cells.Item(3, 11).Value = 5001750
$ws.Cells.Item(3, 13).Value = -5001637
$ws.Cells.Item(16, 8).Value = 2737.625
$ws.Cells.Item(16, 9).Value = 2670.3333
$ws.Cells.Item(16, 10).Value = 2778
$ws.Cells.Item(16, 11).Value = 2670.3333
$ws.Cells.Item(16, 12).Value = 2778
$ws.Cells.Item(16, 13).Value = -2383.3333
$ws.Cells.Item(16, 14).Value = -3352
$ws.Cells.Item(58, 8).Value = 1274.25
$ws.Cells.Item(58, 10).Value = 548.5
$ws.Cells.Item(58, 12).Value = 548.5
$ws.Cells.Item(58, 14).Value = -954.5
$ws.Cells.Item(86, 8).Value = 4498.3335
$ws.Cells.Item(86, 9).Value = 4498.3335
$ws.Cells.Item(86, 11).Value = 4498.3335
$ws.Cells.Item(86, 13).Value = -3375.3335
$ws.Cells.Item(89, 8).Value = 4498.3335
$ws.Cells.Item(89, 9).Value = 4498.3335
$ws.Cells.Item(89, 11).Value = 22491.6675
$ws.Cells.Item(89, 13).Value = -16875.6675
$ws.Cells.Item(100, 8).Value = 99999.5
$ws.Cells.Item(100, 10).Value = 99999.5
$ws.Cells.Item(100, 12).Value = 99999.5
$ws.Cells.Item(100, 14).Value = -102163.5
$ws.Cells.Item(105, 8).Value = 959.0769
$ws.Cells.Item(105, 9).Value = 949.1739
$ws.Cells.Item(105, 11).Value = 949.1739
$ws.Cells.Item(105, 13).Value = 797.8261
$ws.Cells.Item(107, 8).Value = 579.125
$ws.Cells.Item(107, 9).Value = 547.1667
$ws.Cells.Item(107, 11).Value = 547.1667
$ws.Cells.Item(107, 13).Value = 1372.8333
$ws.Cells.Item(113, 8).Value = 2737.625
$ws.Cells.Item(113, 9).Value = 2670.3333
$ws.Cells.Item(113, 10).Value = 2778
$ws.Cells.Item(113, 11).Value = 2670.3333
$ws.Cells.Item(113, 12).Value = 2778
$ws.Cells.Item(113, 13).Value = -500.3332999999998
$ws.Cells.Item(113, 14).Value = -7118
$ws.Cells.Item(136, 8).Value = 1274.25
$ws.Cells.Item(136, 10).Value = 548.5
$ws.Cells.Item(136, 12).Value = 1645.5
$ws.Cells.Item(136, 14).Value = -6745.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1166
$ws.Cells.Item(131, 9).Value = 999
$ws.Cells.Item(131, 10).Value = 1500
$ws.Cells.Item(131, 11).Value = 2997
$ws.Cells.Item(131, 12).Value = 4500
$ws.Cells.Item(131, 13).Value = 2043
$ws.Cells.Item(131, 14).Value = -14580

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 3475.625
$ws.Cells.Item(122, 9).Value = 2232
$ws.Cells.Item(122, 11).Value = 6696
$ws.Cells.Item(122, 13).Value = -4246
$ws.Cells.Item(132, 8).Value = 57954.773
$ws.Cells.Item(132, 9).Value = 72353.586
$ws.Cells.Item(132, 11).Value = 217060.758
$ws.Cells.Item(132, 13).Value = -214530.758

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5083.067
$ws.Cells.Item(122, 10).Value = 5107.143
$ws.Cells.Item(122, 12).Value = 15321.429
$ws.Cells.Item(122, 14).Value = -20221.429
$ws.Cells.Item(136, 8).Value = 8000
$ws.Cells.Item(136, 9).Value = 6000
$ws.Cells.Item(136, 10).Value = 10000
$ws.Cells.Item(136, 11).Value = 18000
$ws.Cells.Item(136, 12).Value = 30000
$ws.Cells.Item(136, 13).Value = -15450
$ws.Cells.Item(136, 14).Value = -35100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 57004
$ws.Cells.Item(3, 10).Value = 57004
$ws.Cells.Item(3, 12).Value = 57004
$ws.Cells.Item(3, 14).Value = -57232
$ws.Cells.Item(4, 8).Value = 7810
$ws.Cells.Item(4, 9).Value = 8879.166999999999
$ws.Cells.Item(4, 10).Value = 3533.3333
$ws.Cells.Item(4, 11).Value = 8879.166999999999
$ws.Cells.Item(4, 12).Value = 3533.3333
$ws.Cells.Item(4, 13).Value = -8766.166999999999
$ws.Cells.Item(4, 14).Value = -3759.3333
$ws.Cells.Item(6, 8).Value = 872.5
$ws.Cells.Item(6, 10).Value = 1500
$ws.Cells.Item(6, 12).Value = 1500
$ws.Cells.Item(6, 14).Value = -1730
$ws.Cells.Item(113, 8).Value = 525.85
$ws.Cells.Item(113, 9).Value = 478.6154
$ws.Cells.Item(113, 11).Value = 1435.8462
$ws.Cells.Item(113, 13).Value = 734.1538
$ws.Cells.Item(126, 8).Value = 4559.5
$ws.Cells.Item(126, 9).Value = 3145.2
$ws.Cells.Item(126, 11).Value = 9435.599999999999
$ws.Cells.Item(126, 13).Value = -6965.599999999999
